# NW 2008.xlsx - "Add cantrals by cantons"
# Adds idx/idx2/Name/Date Start/Date End header columns (data already
# existed in A:E, it just lacked a header row), renames/split the
# (MW)/(GWh) headers into (MW1)/(MW2)/(GWh) Winter/Summer/Year, and
# removes the old Hiver/Eté/Année/units row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "Hiver / Eté / Année / units" row (row 2) - the header
# row directly above the data will be rewritten from scratch below.
$ws.Rows.Item(2).Delete()

# Rewrite the header row (row 1) completely.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 are brand new header cells - make sure they use the plain
# default style (E1 previously held a styled "(m3/s)" header).
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A1:E1").Font.Size = 10
$ws.Range("A1:E1").Font.Name = "Arial"

# Header cells F1:K1 keep the small (9pt) Arial font used elsewhere in
# the sheet (A1:E1 stay in the default font, matching the diff).
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Put the selection where the author's saved view left it.
[void]$ws.Range("A2:K2").Select()
